# TC09_C3DC_phs002371_TrtmntType-Chemotherapy.xlsx
# Updated C3DC phs002371 queries:
#  - Treatment query (row 5 / cell B5) gets an extra filter
#    "AND trt.treatment_id IS NOT NULL" appended to its WHERE clause.
#  - The Treatment Response (row 6) and Survival (row 7) queries are
#    unchanged in content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$treatmentQuery = @"
SELECT
    DISTINCT prt.participant_id AS "Participant Id",
    trt.treatment_id AS "Treatment Id",
    CASE 
    WHEN trt.age_at_treatment_start = -999 THEN 'Not Reported'
    WHEN trt.age_at_treatment_start >= 1000 THEN 
        substr(trt.age_at_treatment_start, 1, length(trt.age_at_treatment_start) - 3) || ',' || substr(trt.age_at_treatment_start, -3)
    ELSE 
        trt.age_at_treatment_start 
END AS "Age at Treatment Start",
    CASE 
    WHEN trt.age_at_treatment_end = -999 THEN 'Not Reported'
    WHEN trt.age_at_treatment_end >= 1000 THEN 
        substr(trt.age_at_treatment_end, 1, length(trt.age_at_treatment_end) - 3) || ',' || substr(trt.age_at_treatment_end, -3)
    ELSE 
        trt.age_at_treatment_end 
END AS "Age at Treatment End",
    trt.treatment_type AS "Treatment Type",
    REPLACE(trt.treatment_agent, ';', ', ') AS "Treatment Agent",
    std.dbgap_accession AS "dbGaP Accession"
FROM 
    df_study std
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    df_diagnoses dgn ON prt.id = dgn."participant.id"
LEFT JOIN 
    df_treatments trt ON prt.id = trt."participant.id"
LEFT JOIN 
    df_treatment_resp trr ON prt.id = trr."participant.id"
LEFT JOIN 
    df_survival srv ON prt.id = srv."participant.id"
LEFT JOIN 
    df_reference_files rfs ON std.id = rfs."study.id"
WHERE 
    std.dbgap_accession = 'phs002371' AND trt.treatment_type = 'Chemotherapy' AND trt.treatment_id IS NOT NULL
ORDER BY 
    trt.treatment_id ASC
LIMIT 100;
"@

# Row 5 (TreatmentTab): rewrite the query text and nudge the font size
# down to 11pt (it was 12pt, matching the other query cells).
$b5 = $ws.Range("B5")
$b5.Value = $treatmentQuery
$b5.Font.Size = 11

# Leave rows 6 (TreatmentRespTab) / 7 (SurvivalTab) content untouched -
# only the selection/scroll position moved as a side effect of editing
# row 5 in the real workbook.
$ws.Range("C5").Select() | Out-Null

Write-Host "B5 length:" $ws.Range("B5").Text.Length
Write-Host "B5 font size:" $ws.Range("B5").Font.Size
Write-Host "B6 font size:" $ws.Range("B6").Font.Size
Write-Host "B7 font size:" $ws.Range("B7").Font.Size
